$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = -196771045.85
$ws.Range("P2").Value = -1900.4379758512
$ws.Range("Q2").Value = 1010218615.95
$ws.Range("R2").Value = 9756.8105783551
$ws.Range("S2").Value = 39946943.73
$ws.Range("T2").Value = 385.8122954815
$ws.Range("U2").Value = -10935810.88
$ws.Range("V2").Value = -105.6193516851
$ws.Range("Y2").Value = 31659766.42
$ws.Range("Z2").Value = 305.773759301
$ws.Range("AA2").Value = 198772977.35
$ws.Range("AB2").Value = 1919.772866466
$ws.Range("AC2").Value = -10353984.1
$ws.Range("AD2").Value = -353.2882313242
